$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header): add P1=14, Q1=15, matching style of existing header cells (s="1")
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Rows 2-25: update columns I, K, M, O and add columns P, Q
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2 (new)
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2 (new)
}
